# Updates for "Update code for caclular Luong ca nhan"
#
# 1. Insert a new worksheet "Đơn phụ phẫu 1" between "Đơn sale chính" and
#    "Lương", populated with the per-order detail (header + 1 data row +
#    total row).
# 2. Update the "Lương" sheet: new "Ngày công"/"Phụ cấp" totals, a new
#    "Ứng lương" line per cơ sở, renumbered rows and new "Tổng lương" rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New sheet "Đơn phụ phẫu 1"
# ---------------------------------------------------------------------
$luong = $wb.Worksheets.Item("Lương")
$donPhuPhau = $wb.Worksheets.Add($luong)
$donPhuPhau.Name = "Đơn phụ phẫu 1"

# Re-resolve the "Lương" handle: Worksheets.Add() shifted what the old
# $luong reference pointed at (it now resolves to the new sheet instead
# of following the original one), so fetch it again by name.
$luong = $wb.Worksheets.Item("Lương")

$headers = @(
    "Tiền tố", "Mã dịch vụ", "Ngày thực hiện", "Cơ sở", "Khách hàng",
    "Nguồn khách", "Nhóm dịch vụ", "Tên dịch vụ", "Sale chính",
    "Đơn giá gốc", "Sale phụ", "Upsale", "Đơn giá", "Thanh toán lần đầu",
    "Trả sau", "Đã thanh toán", "Dư nợ", "Bác sĩ 1", "Bác sĩ 2",
    "Phụ phẫu 1", "Phụ phẫu 2", "Công phụ phẫu 1", "Công phụ phẫu 2",
    "Tỉ lệ chiết khấu sale chính", "Tỉ lệ chiết khấu sale phụ",
    "Chiết khấu sale chính", "Chiết khấu sale phụ"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $donPhuPhau.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$row2 = @{
    1 = "HD-LUXURY"; 2 = 554; 3 = "07-15-2024"; 4 = "SÓC TRĂNG";
    5 = "lê thị bích"; 6 = "Cá nhân"; 7 = "Đại phẫu"; 8 = "cắt sẹo ";
    9 = "Lâm Thị Mỹ Hằng"; 10 = 7000000; 11 = "Lê Đình Hậu"; 12 = 6000000;
    13 = 13000000; 14 = 13000000; 15 = 0; 16 = 13000000; 17 = 0;
    18 = "Bác Sĩ Thảo"; 20 = "Trần Khánh Hiệp";
    24 = 0.1; 25 = 0.04; 26 = 1060000; 27 = 240000
}
# Column 3 ("Ngày thực hiện") holds a date-look-alike string that must
# stay literal text instead of being auto-parsed into a date serial.
$donPhuPhau.Cells.Item(2, 3).NumberFormat = "@"
foreach ($col in $row2.Keys) {
    $donPhuPhau.Cells.Item(2, $col).Value = $row2[$col]
}

$row3 = @{
    1 = "Tổng"; 2 = 1; 10 = 7000000; 12 = 6000000; 13 = 13000000;
    14 = 13000000; 15 = 0; 16 = 13000000; 17 = 0; 22 = 0; 23 = 0;
    24 = 0.1; 25 = 0.04; 26 = 1060000; 27 = 240000
}
foreach ($col in $row3.Keys) {
    $donPhuPhau.Cells.Item(3, $col).Value = $row3[$col]
}

# ---------------------------------------------------------------------
# 2. Updated "Lương" sheet
# ---------------------------------------------------------------------
$luongRows = @(
    @("Danh mục", 8),
    @("Ngày công", 14),
    @("Phụ cấp", 490000),
    @("Lương cơ bản tại CẦN THƠ", $null),
    @("Chiết khấu sale chính tại CẦN THƠ", 0),
    @("Chiết khấu sale phụ tại CẦN THƠ", 0),
    @("Đơn 1 bác sĩ tại CẦN THƠ", 0),
    @("Đơn 2 bác sĩ tại CẦN THƠ", 0),
    @("Công phụ phẫu 1 tại CẦN THƠ", 0),
    @("Công phụ phẫu 2 tại CẦN THƠ", 0),
    @("Ứng lương tại CẦN THƠ", 0),
    @("Lương cơ bản tại LONG XUYÊN", $null),
    @("Chiết khấu sale chính tại LONG XUYÊN", 0),
    @("Chiết khấu sale phụ tại LONG XUYÊN", 0),
    @("Đơn 1 bác sĩ tại LONG XUYÊN", 0),
    @("Đơn 2 bác sĩ tại LONG XUYÊN", 0),
    @("Công phụ phẫu 1 tại LONG XUYÊN", 0),
    @("Công phụ phẫu 2 tại LONG XUYÊN", 0),
    @("Ứng lương tại LONG XUYÊN", 0),
    @("Lương cơ bản tại SÓC TRĂNG", 2060000),
    @("Chiết khấu sale chính tại SÓC TRĂNG", 0),
    @("Chiết khấu sale phụ tại SÓC TRĂNG", 0),
    @("Đơn 1 bác sĩ tại SÓC TRĂNG", 0),
    @("Đơn 2 bác sĩ tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 1 tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 2 tại SÓC TRĂNG", 0),
    @("Ứng lương tại SÓC TRĂNG", -1000000),
    @("Tổng lương tại CẦN THƠ", 0),
    @("Tổng lương tại LONG XUYÊN", 0),
    @("Tổng lương tại SÓC TRĂNG", 5225000),
    @("Tổng lương", 5225000)
)

for ($i = 0; $i -lt $luongRows.Length; $i++) {
    $r = $i + 1
    $luong.Cells.Item($r, 1).Value = $luongRows[$i][0]
    $val = $luongRows[$i][1]
    if ($null -ne $val) {
        $luong.Cells.Item($r, 2).Value = $val
    } else {
        $luong.Cells.Item($r, 2).Value = ""
    }
}
